# Generate Report for handoff
# Update the "Latest Handoff Datetime" (column D) for the b96b9ac6 file row (row 5)
# on both the zh-cn and de-de worksheets, reflecting a fresh handoff report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-18 03:19:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-18 03:19:48"
